$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the existing header row (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-47
$data = @{
    2  = @(3,3)
    3  = @(8,9)
    4  = @(9,9)
    5  = @(8,8)
    6  = @(2,3)
    7  = @(8,8)
    8  = @(9,9)
    9  = @(7,8)
    10 = @(9,9)
    11 = @(7,8)
    12 = @(6,6)
    13 = @(9,9)
    14 = @(8,8)
    15 = @(4,5)
    16 = @(5,5)
    17 = @(8,8)
    18 = @(7,7)
    19 = @(6,6)
    20 = @(5,6)
    21 = @(9,9)
    22 = @(5,5)
    23 = @(6,7)
    24 = @(8,8)
    25 = @(6,7)
    26 = @(7,8)
    27 = @(7,7)
    28 = @(8,8)
    29 = @(8,9)
    30 = @(7,7)
    31 = @(4,5)
    32 = @(6,7)
    33 = @(4,6)
    34 = @(6,7)
    35 = @(6,7)
    36 = @(7,7)
    37 = @(8,9)
    38 = @(7,7)
    39 = @(7,8)
    40 = @(6,6)
    41 = @(2,2)
    42 = @(4,4)
    43 = @(8,8)
    44 = @(4,5)
    45 = @(5,5)
    46 = @(8,8)
    47 = @(8,8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

# Update the sheet dimension to reflect the new used range
$ws.UsedRange | Out-Null
